# Fill in the TestCases & resultaten sheet with the system-test results for
# rows 3 and 4 (Volgnr 1 and 2), then update the active sheet / selection
# so TestCases & resultaten becomes the visible tab (matching the saved
# workbook view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases & resultaten")

# --- Row 3 (Volgnr 1) ---
$ws.Range("B3").Value = 'Elementor'
$ws.Range("C3").Value = 'Tot nu toe alleen het formulier toegevoegd in Elementor, maar de pagina''s nog niks mee gedaan vanwege er te weinig tijd was'
$ws.Range("D3").Value = 'Inloggen in Elementor en kijken of er iets is gemaakt.'
$ws.Range("E3").Value = 'Niet aan begonnen, niks gemaakt.'
$ws.Range("F3").Value = 'Ja'
$ws.Range("G3").Value = 'Dat er geen één pagina is gemaakt. Aangezien er geen tijd meer hiervoor was. '
$ws.Range("H3").Value = 'Gekeken naar het gehele project in Elementor, alleen de aanpassingen die er zijn gemaakt is de kleuren en logo die erin zitten. '
$ws.Range("I3").Value = 'Ja'
$ws.Range("J3").Value = 'Laag'
$ws.Range("K3").Value = 'Philip Klok'
$ws.Range("L3").Value = "5/5/2023"
$ws.Range("M3").Value = 'Na de PvB wordt hier naar gewerkt en denkt aan het spreken van andere medewerkers in het project. '
$ws.Range("N3").Value = 'Wordt nog vervolgd'
$ws.Range("O3").Value = 'Philip Klok'
$ws.Range("P3").Value = 'Wordt nog vervolgd'
$ws.Range("Q3").Value = 'Nee'

# --- Row 4 (Volgnr 2) ---
$ws.Range("B4").Value = 'Multi-step form'
$ws.Range("C4").Value = 'Het formulier is af en zijn in gesprek gegaan over het formulier, volgens Pippijn Stortelder is het formulier compleet.'
$ws.Range("D4").Value = 'Of er nog foutmeldingen zijn en errors in de code/ Elementor'
$ws.Range("E4").Value = 'Compleet, niks meer eraan doen.'
$ws.Range("F4").Value = 'Nee'
$ws.Range("G4").Value = 'In de code mist nog wat commentaar en in Elementor staat bij het resultaat de tekst niet in het midden. '
$ws.Range("H4").Value = 'Eerst naar de code wezen kijken en de opmaak zag er goed uit, de code was mooi verdeeld, alleen wat het nog wat misten was het commentaar bij het JavaScript gedeelte. In Elementor keken we naar het formulier en gingen alle stappen af, ook weer terug bij stap 2 om te kijken of het werkten.'
$ws.Range("I4").Value = 'Ja'
$ws.Range("J4").Value = 'Hoog'
$ws.Range("K4").Value = 'Philip Klok'
$ws.Range("L4").Value = "5/8/2023"
$ws.Range("M4").Value = 'Bij het commentaar in de code, zorg ervoor dat je het commentaar goed mogelijk uitleg zodat andere de functies snappen.'
$ws.Range("N4").Value = "5/8/2023"
$ws.Range("O4").Value = 'Philip Klok'
$ws.Range("P4").Value = 'De lijnen van code opnieuw naar gekeken en daarna uitgelegd boven de functie.'
$ws.Range("Q4").Value = 'Nee'

# --- Window / selection state ---
# TestPlan's stored selection moves from A2 to C3.
$wsPlan = $wb.Worksheets.Item("TestPlan")
$wsPlan.Range("C3").Select()

# TestCases & resultaten becomes the active (visible) tab, with the
# active cell resting on J7 inside the frozen bottom-right pane.
$ws.Activate()
$ws.Range("J7").Select()
